# Update "想去人数" (F column) values for specific events on the
# 展览 (Exhibition) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4088
$ws1.Range("F11").Value = 296
$ws1.Range("F13").Value = 2868
$ws1.Range("F15").Value = 1284

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 4088
$ws4.Range("F12").Value = 296
$ws4.Range("F14").Value = 2868
$ws4.Range("F16").Value = 1284
